$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows with new word/meaning pairs
$ws.Range("A2").Value = "('Word', 'moutaine')"
$ws.Range("B2").Value = "('Meaning', 'son tinh')"

$ws.Range("A3").Value = "('Word', 'good bye')"
$ws.Range("B3").Value = "('Meaning', 'tam biet')"

# Remove the now-obsolete rows 4-6, shifting remaining rows up
$ws.Range("A4:B6").EntireRow.Delete()
